# Fruta / hortaliza, semanal
# Insert two new weekly price rows (for 2022-11-30) above the existing
# "Feria Lagunitas de Puerto Montt - Arándano (blue)" record that used to
# sit at row 18, pushing the rest of the table down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at row 18; this shifts the old rows 18-24 down to 20-26
# and keeps the existing column D date-style formatting on the new rows.
$ws.Rows("18:19").Insert()

# --- New row 18 ---
$ws.Cells.Item(18, 1).Value = 4
$ws.Cells.Item(18, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(18, 3).Value = "Los Lagos"
$ws.Cells.Item(18, 4).Value = 44895
$ws.Cells.Item(18, 5).Value = 10
$ws.Cells.Item(18, 6).Value = "Fruta"
$ws.Cells.Item(18, 7).Value = 100101
$ws.Cells.Item(18, 8).Value = "Berries"
$ws.Cells.Item(18, 9).Value = 100101001
$ws.Cells.Item(18, 10).Value = "Arándano (blue)"
$ws.Cells.Item(18, 11).Value = "Sin especificar"
$ws.Cells.Item(18, 12).Value = "Primera"
$ws.Cells.Item(18, 13).Value = 120
$ws.Cells.Item(18, 14).Value = 8000
$ws.Cells.Item(18, 15).Value = 8500
$ws.Cells.Item(18, 16).Value = 8250
$ws.Cells.Item(18, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(18, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(18, 19).Value = 4125
$ws.Cells.Item(18, 20).Value = 2

# --- New row 19 ---
$ws.Cells.Item(19, 1).Value = 4
$ws.Cells.Item(19, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(19, 3).Value = "Los Lagos"
$ws.Cells.Item(19, 4).Value = 44895
$ws.Cells.Item(19, 5).Value = 10
$ws.Cells.Item(19, 6).Value = "Fruta"
$ws.Cells.Item(19, 7).Value = 100101
$ws.Cells.Item(19, 8).Value = "Berries"
$ws.Cells.Item(19, 9).Value = 100101001
$ws.Cells.Item(19, 10).Value = "Arándano (blue)"
$ws.Cells.Item(19, 11).Value = "Sin especificar"
$ws.Cells.Item(19, 12).Value = "Primera"
$ws.Cells.Item(19, 13).Value = 300
$ws.Cells.Item(19, 14).Value = 8000
$ws.Cells.Item(19, 15).Value = 8500
$ws.Cells.Item(19, 16).Value = 8250
$ws.Cells.Item(19, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(19, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(19, 19).Value = 4125
$ws.Cells.Item(19, 20).Value = 2
